$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12213.0524555979
$ws.Range("C2").Value = 11536.8010635582
$ws.Range("E2").Value = 7802.70116365903
$ws.Range("F2").Value = 46.6517594673843

$ws.Range("B3").Value = 11728.2894823166
$ws.Range("C3").Value = 11468.0431508643
$ws.Range("E3").Value = 7419.77098525419
$ws.Range("F3").Value = 283.831422338269

$ws.Range("B4").Value = 11900.4253798072
$ws.Range("C4").Value = 11030.9537607406
$ws.Range("E4").Value = 8018.17889670135
$ws.Range("F4").Value = 290.553027393416

$ws.Range("B5").Value = 4807.81477322236
$ws.Range("C5").Value = 7929.50299605596
$ws.Range("E5").Value = 7910.54191568641
$ws.Range("F5").Value = 156.841037989265

$ws.Range("B6").Value = 5266.01744988633
$ws.Range("C6").Value = 7923.73630309195
$ws.Range("E6").Value = 7949.00893782206
$ws.Range("F6").Value = 253.370218371417

$ws.Range("B7").Value = 11794.4383606187
$ws.Range("C7").Value = 11093.5894274355
$ws.Range("E7").Value = 7860.08675015067
$ws.Range("F7").Value = 381.742340732755

$ws.Range("B8").Value = 11503.27393928
$ws.Range("C8").Value = 10618.8381073625
$ws.Range("E8").Value = 7642.4222902102
$ws.Range("F8").Value = 352.891683232194

$ws.Range("B9").Value = 11503.27393928
$ws.Range("C9").Value = 10531.4904497926
$ws.Range("E9").Value = 7642.4222902102
$ws.Range("F9").Value = 349.252197500117

$ws.Range("B10").Value = 11503.27393928
$ws.Range("C10").Value = 10683.4908512109
$ws.Range("E10").Value = 7642.4222902102
$ws.Range("F10").Value = 355.585547559213

$ws.Range("B11").Value = 11503.27393928
$ws.Range("C11").Value = 10022.5901213771
$ws.Range("E11").Value = 7642.4222902102
$ws.Range("F11").Value = 328.048017149469

$ws.Range("B12").Value = 4875.64205484986
$ws.Range("C12").Value = 7197.1206598242
$ws.Range("E12").Value = 7273.63201759004
$ws.Range("F12").Value = 194.953861558927

$ws.Range("B13").Value = 4779.1830818516
$ws.Range("C13").Value = 7242.28414327488
$ws.Range("E13").Value = 7266.27291865256
$ws.Range("F13").Value = 196.529044246977

$ws.Range("B14").Value = 11115.7963561476
$ws.Range("C14").Value = 10513.1518931632
$ws.Range("E14").Value = 7172.64234095707
$ws.Range("F14").Value = 328.913926421676

$ws.Range("B15").Value = 11115.7963561476
$ws.Range("C15").Value = 10776.8868848665
$ws.Range("E15").Value = 7172.64234095707
$ws.Range("F15").Value = 339.902884409315
